# Add six new live-cam rows (265-270) to the "location-1" sheet, each with a
# Category, Lat/Long, Location name, City/State, Country, YouTube video id and
# a YouTube hyperlink in column G (mirroring the existing rows 261-264).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$hyperlinkUrl = "https://www.youtube.com/@BostonAndMaineLive/streams"
$hyperlinkText = "(170) Boston and Maine Live - YouTube"

$rows = @(
    @{ A = "LIVE, CITY, BUILDING, METRO";  B = "42.36157044785626, -71.07530945189643"; C = "Boston, MA Live Cam - City of Boston"; D = "MA"; E = "USA";    F = "cSCXEW0tda8" },
    @{ A = "LIVE, SEA, CITY, BUILDING";    B = "42.35796880969909, -71.04723638744522"; C = "Boston Harbor, Massachusetts - Live"; D = "MA"; E = "USA";    F = "wRBSMzhZL50" },
    @{ A = "LIVE, RAIL, TRAIN, STATION";   B = "45.59023541739565, -67.32878944964192"; C = "McAdam Railway Station - Live Cam";   D = "NB"; E = "Canada"; F = "P8t1nuM8BcY" },
    @{ A = "LIVE, MOUNTAIN, SKI, RESORT";  B = "44.056082290650515, -71.6275010341407"; C = "Loon Mountain, New Hampshire - The Mountain Club at Loon - Live Cam"; D = "NH"; E = "USA"; F = "2l-EzK0erOE" },
    @{ A = "LIVE, DAM, RIVER";             B = "45.183980846203625, -69.23010764376563"; C = "Mayo Mill Dam at Dover-Foxcroft, Maine LIVE cam"; D = "ME"; E = "USA"; F = "JK9D1UPy6s0" },
    @{ A = "LIVE, MOUNTAIN, NATURE";       B = "43.08084130296787, 11.711669107435554"; C = "Tuscany, Italy Live Webcam - Podere Il Casale"; D = "Pienza"; E = "Italy"; F = "vnDRqc0GCaI" }
)

$startRow = 265
$templateRow = 264

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]

    # Create the hyperlink first (this also sets G's value/display to the URL);
    # we overwrite the cell text afterwards but keep the "display" = URL, just
    # like the existing rows in the sheet.
    $ws.Hyperlinks.Add($ws.Range("G$r"), $hyperlinkUrl, "", "", $hyperlinkUrl) | Out-Null

    # Clone the direct formatting (borders / hyperlink style) used by the
    # template row so the new rows look consistent with the existing table.
    $ws.Range("A$templateRow").Copy() | Out-Null
    $ws.Range("A$r").PasteSpecial(-4122) | Out-Null
    $ws.Range("E$templateRow").Copy() | Out-Null
    $ws.Range("E$r").PasteSpecial(-4122) | Out-Null
    $ws.Range("G$templateRow").Copy() | Out-Null
    $ws.Range("G$r").PasteSpecial(-4122) | Out-Null

    $ws.Cells.Item($r, 1).Value2 = $row.A
    $ws.Cells.Item($r, 2).Value2 = $row.B
    $ws.Cells.Item($r, 3).Value2 = $row.C
    $ws.Cells.Item($r, 4).Value2 = $row.D
    $ws.Cells.Item($r, 5).Value2 = $row.E
    $ws.Cells.Item($r, 6).Value2 = $row.F
    $ws.Cells.Item($r, 7).Value2 = $hyperlinkText
}

# Keep the header row frozen and scroll/select the same way the source
# workbook ends up after the new rows are appended.
$ws.Activate()
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$excel.Goto($ws.Range("A244"))
$ws.Range("A271").Select() | Out-Null
